# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G, header "K" in row 1) previously held the
# Strike# values for each start; this regenerates that column with the
# actual strikeout (K) counts for each of the 34 starts (rows 2-35).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 2
    6  = 1
    7  = 0
    8  = 3
    9  = 1
    10 = 5
    11 = 4
    12 = 1
    13 = 5
    14 = 5
    15 = 7
    16 = 3
    17 = 7
    18 = 4
    19 = 4
    20 = 4
    21 = 7
    22 = 3
    23 = 3
    24 = 6
    25 = 4
    26 = 5
    27 = 7
    28 = 3
    29 = 1
    30 = 5
    31 = 4
    32 = 6
    33 = 3
    34 = 3
    35 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
